$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell B1: make it look like A1 (bold, bordered, centered, integer format) ---
$ws.Range("B1").NumberFormat = $ws.Range("A1").NumberFormat
$ws.Range("B1").Font.Bold = $ws.Range("A1").Font.Bold
$ws.Range("B1").Borders.LineStyle = $ws.Range("A1").Borders.LineStyle
$ws.Range("B1").HorizontalAlignment = $ws.Range("A1").HorizontalAlignment
$ws.Range("B1").VerticalAlignment = $ws.Range("A1").VerticalAlignment

# --- Column B body cells (B2:B12): switch from 2-decimal text-backed prices to plain integer numbers ---
$ws.Range("B2:B12").NumberFormat = "0"

# Replace the price values that were stored as shared-string text with real numeric values.
$ws.Range("B2").Value = 41800
$ws.Range("B3").Value = 44000
$ws.Range("B4").Value = 37000
$ws.Range("B5").Value = 31200
$ws.Range("B6").Value = 333333
$ws.Range("B7").Value = 31000
$ws.Range("B8").Value = 29655
$ws.Range("B9").Value = 31105
$ws.Range("B10").Value = 31105
$ws.Range("B11").Value = 37699
$ws.Range("B12").Value = 37699

# --- Update the remembered selection to match the new edit location ---
$ws.Range("I14").Select()
